$wb = $excel.ActiveWorkbook

# Add the new worksheet and rename it to "promotion"
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "promotion"

# Move it to the end of the tab order (after the current last sheet, "Sheet1")
$newSheet.Move($null, $wb.Worksheets.Item("Sheet1"))

# Re-acquire the worksheet reference by name since the Move() call can
# reseat/invalidate the old reference
$ws = $wb.Worksheets.Item("promotion")

# Header row
$ws.Range("A1").Value = "t1.A"
$ws.Range("B1").Value = "t1.B"
$ws.Range("C1").Value = "t2.A"
$ws.Range("D1").Value = "t2.B"
$ws.Range("E1").Value = "t3.A"
$ws.Range("F1").Value = "t3.B"
$ws.Range("G1").Value = "t3.C"

# Data row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.5
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = 0.5

# Select G2 and make "promotion" the active sheet/tab
$ws.Range("G2").Select()
$ws.Activate()
